$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats for the id column (A) and date column (E) from the last existing row (161)
$ws.Cells.Item(161, 1).Copy() | Out-Null
$ws.Cells.Item(162, 1).Resize(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(161, 5).Copy() | Out-Null
$ws.Cells.Item(162, 5).Resize(5, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 162 (id=160)
$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(162, 2).Value = 7939469
$ws.Cells.Item(162, 3).Value = "Hungary NB I"
$ws.Cells.Item(162, 4).Value = "Hungary NB I"
$ws.Cells.Item(162, 5).Value = 45392.54166666666
$ws.Cells.Item(162, 6).Value = "Ferencvarosi TC"
$ws.Cells.Item(162, 7).Value = "Diosgyori VTK"
$ws.Cells.Item(162, 11).Value = 1.3
$ws.Cells.Item(162, 12).Value = 5.5
$ws.Cells.Item(162, 13).Value = 9
$ws.Cells.Item(162, 14).Value = 1.3
$ws.Cells.Item(162, 15).Value = 5.75
$ws.Cells.Item(162, 16).Value = 8.5
$ws.Cells.Item(162, 17).Value = -1.5
$ws.Cells.Item(162, 18).Value = 1.825
$ws.Cells.Item(162, 19).Value = 2.025
$ws.Cells.Item(162, 20).Value = 3.25
$ws.Cells.Item(162, 21).Value = 2.05
$ws.Cells.Item(162, 22).Value = 1.8
$ws.Cells.Item(162, 23).Value = 0
$ws.Cells.Item(162, 24).Value = 0
$ws.Cells.Item(162, 25).Value = 0
$ws.Cells.Item(162, 26).Value = 0
$ws.Cells.Item(162, 27).Value = 0

# Row 163 (id=161)
$ws.Cells.Item(163, 1).Value = 161
$ws.Cells.Item(163, 2).Value = 6818359
$ws.Cells.Item(163, 3).Value = "Hungary NB I"
$ws.Cells.Item(163, 4).Value = "Hungary NB I"
$ws.Cells.Item(163, 5).Value = 45395.39583333334
$ws.Cells.Item(163, 6).Value = "Kisvarda FC"
$ws.Cells.Item(163, 7).Value = "Ujpest"
$ws.Cells.Item(163, 11).Value = 2.6
$ws.Cells.Item(163, 12).Value = 3.4
$ws.Cells.Item(163, 13).Value = 2.6
$ws.Cells.Item(163, 14).Value = 2.6
$ws.Cells.Item(163, 15).Value = 3.4
$ws.Cells.Item(163, 16).Value = 2.6
$ws.Cells.Item(163, 17).Value = 0
$ws.Cells.Item(163, 18).Value = 1.925
$ws.Cells.Item(163, 19).Value = 1.925
$ws.Cells.Item(163, 20).Value = 2.5
$ws.Cells.Item(163, 21).Value = 1.975
$ws.Cells.Item(163, 22).Value = 1.875
$ws.Cells.Item(163, 23).Value = 0
$ws.Cells.Item(163, 24).Value = 0
$ws.Cells.Item(163, 25).Value = 0
$ws.Cells.Item(163, 26).Value = 0
$ws.Cells.Item(163, 27).Value = 0

# Row 164 (id=162)
$ws.Cells.Item(164, 1).Value = 162
$ws.Cells.Item(164, 2).Value = 6818358
$ws.Cells.Item(164, 3).Value = "Hungary NB I"
$ws.Cells.Item(164, 4).Value = "Hungary NB I"
$ws.Cells.Item(164, 5).Value = 45395.5
$ws.Cells.Item(164, 6).Value = "Debreceni VSC"
$ws.Cells.Item(164, 7).Value = "Kecskemeti TE"
$ws.Cells.Item(164, 11).Value = 1.666
$ws.Cells.Item(164, 12).Value = 3.5
$ws.Cells.Item(164, 13).Value = 5.5
$ws.Cells.Item(164, 14).Value = 1.727
$ws.Cells.Item(164, 15).Value = 3.4
$ws.Cells.Item(164, 16).Value = 5
$ws.Cells.Item(164, 17).Value = -0.75
$ws.Cells.Item(164, 18).Value = 1.95
$ws.Cells.Item(164, 19).Value = 1.9
$ws.Cells.Item(164, 20).Value = 2.5
$ws.Cells.Item(164, 21).Value = 1.975
$ws.Cells.Item(164, 22).Value = 1.875
$ws.Cells.Item(164, 23).Value = 0
$ws.Cells.Item(164, 24).Value = 0
$ws.Cells.Item(164, 25).Value = 0
$ws.Cells.Item(164, 26).Value = 0
$ws.Cells.Item(164, 27).Value = 0

# Row 165 (id=163)
$ws.Cells.Item(165, 1).Value = 163
$ws.Cells.Item(165, 2).Value = 6818360
$ws.Cells.Item(165, 3).Value = "Hungary NB I"
$ws.Cells.Item(165, 4).Value = "Hungary NB I"
$ws.Cells.Item(165, 5).Value = 45395.60416666666
$ws.Cells.Item(165, 6).Value = "Paksi"
$ws.Cells.Item(165, 7).Value = "MOL Fehervar FC"
$ws.Cells.Item(165, 11).Value = 1.666
$ws.Cells.Item(165, 12).Value = 3.5
$ws.Cells.Item(165, 13).Value = 5.5
$ws.Cells.Item(165, 14).Value = 1.571
$ws.Cells.Item(165, 15).Value = 3.75
$ws.Cells.Item(165, 16).Value = 6
$ws.Cells.Item(165, 17).Value = -1
$ws.Cells.Item(165, 18).Value = 2.05
$ws.Cells.Item(165, 19).Value = 1.8
$ws.Cells.Item(165, 20).Value = 2.75
$ws.Cells.Item(165, 21).Value = 1.85
$ws.Cells.Item(165, 22).Value = 2
$ws.Cells.Item(165, 23).Value = 0
$ws.Cells.Item(165, 24).Value = 0
$ws.Cells.Item(165, 25).Value = 0
$ws.Cells.Item(165, 26).Value = 0
$ws.Cells.Item(165, 27).Value = 0

# Row 166 (id=164)
$ws.Cells.Item(166, 1).Value = 164
$ws.Cells.Item(166, 2).Value = 6818361
$ws.Cells.Item(166, 3).Value = "Hungary NB I"
$ws.Cells.Item(166, 4).Value = "Hungary NB I"
$ws.Cells.Item(166, 5).Value = 45396.57291666666
$ws.Cells.Item(166, 6).Value = "Mezokovesd Zsory"
$ws.Cells.Item(166, 7).Value = "Puskas Academy"
$ws.Cells.Item(166, 11).Value = 5.5
$ws.Cells.Item(166, 12).Value = 3.5
$ws.Cells.Item(166, 13).Value = 1.666
$ws.Cells.Item(166, 14).Value = 5.25
$ws.Cells.Item(166, 15).Value = 3.5
$ws.Cells.Item(166, 16).Value = 1.7
$ws.Cells.Item(166, 17).Value = 0.75
$ws.Cells.Item(166, 18).Value = 1.925
$ws.Cells.Item(166, 19).Value = 1.925
$ws.Cells.Item(166, 20).Value = 2.5
$ws.Cells.Item(166, 21).Value = 1.925
$ws.Cells.Item(166, 22).Value = 1.925
$ws.Cells.Item(166, 23).Value = 0
$ws.Cells.Item(166, 24).Value = 0
$ws.Cells.Item(166, 25).Value = 0
$ws.Cells.Item(166, 26).Value = 0
$ws.Cells.Item(166, 27).Value = 0
